$d = $word.ActiveDocument

function Insert-XmlFragment($range, [string]$bodyXml) {
    $pkg = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData>' +
           '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
           '<w:body>' + $bodyXml + '</w:body>' +
           '</w:document>' +
           '</pkg:xmlData></pkg:part></pkg:package>'
    $range.InsertXML($pkg)
}

# --- 1. "Esta es otra frase para probar el commit" -> split "commit" into its
#        own run, wrapped in proofErr spellStart/spellEnd markers. ---
$pCommit = $d.Paragraphs(4)
$rCommit = $pCommit.Range
$commitXml = '<w:p><w:r><w:t xml:space="preserve">Esta es otra frase para probar el </w:t></w:r>' +
             '<w:proofErr w:type="spellStart"/><w:r><w:t>commit</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>'
Insert-XmlFragment $rCommit $commitXml

# --- 2. New paragraph "Esta es la frase de prueba creada por Jeiny" + "." ---
#        inserted right after, before the (still underlined) Jeiny paragraph.
$pCommit = $d.Paragraphs(4)
$pCommit.Range.InsertParagraphAfter()
$pJeinyNew = $d.Paragraphs(5)
$jeinyXml = '<w:p><w:r><w:t>Esta es la frase de prueba creada por Jeiny</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p>'
Insert-XmlFragment $pJeinyNew.Range $jeinyXml

# --- 3. The old underlined "Esta es la frase de prueba creada por Jeiny."
#        paragraph loses its run text, staying empty (keeps pPr/underline). ---
$pUnderline = $d.Paragraphs(6)
$rUnderline = $pUnderline.Range
$delRange = $d.Range($rUnderline.Start, $rUnderline.End - 1)
$delRange.Text = ""

# --- 4. Append the new trailing paragraphs after the (now empty) underline
#        paragraph: "Error 1", 4 blanks, "Error raro ", 5 blanks, and a final
#        underlined paragraph with the "Errorrrrrrrr ..." content. ---
$pUnderline = $d.Paragraphs(6)
$anchor = $pUnderline.Range
for ($i = 0; $i -lt 12; $i++) {
    $anchor.InsertParagraphAfter()
    $anchor = $d.Paragraphs(6 + $i + 1).Range
}

$pError1 = $d.Paragraphs(7)
Insert-XmlFragment $pError1.Range '<w:p><w:r><w:t>Error 1</w:t></w:r></w:p>'

$pErrorRaro = $d.Paragraphs(12)
Insert-XmlFragment $pErrorRaro.Range '<w:p><w:r><w:t xml:space="preserve">Error raro </w:t></w:r></w:p>'

$pErrorFinal = $d.Paragraphs(18)
$finalXml = '<w:p><w:pPr><w:rPr><w:u w:val="single"/></w:rPr></w:pPr>' +
            '<w:proofErr w:type="spellStart"/><w:r><w:t>Errorrrrrrrr</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
            '<w:r><w:t xml:space="preserve"> &#161;!!!!&#8221;!&#8221;!!</w:t></w:r>' +
            '<w:r><w:tab/><w:t>&#8220;!</w:t></w:r></w:p>'
Insert-XmlFragment $pErrorFinal.Range $finalXml
